# Update bulk_user_details.xlsx with new iAuthor testcase data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "IjYJr159"
$ws.Range("B2").Value = 23101809
$ws.Range("C2").Value = "kihysnc44"
$ws.Range("D2").Value = "m`$7%y4TU"
$ws.Range("F2").Value = "qLsJbtSe"
$ws.Range("G2").Value = "uLAN"

# Row 3
$ws.Range("A3").Value = "EXPie593"
$ws.Range("B3").Value = 23101808
$ws.Range("C3").Value = "hbczjvr56"
$ws.Range("D3").Value = "jC9c!2%R"
$ws.Range("F3").Value = "VKuUJqdx"
$ws.Range("G3").Value = "EpZa"

# Row 4
$ws.Range("A4").Value = "BabUp443"
$ws.Range("B4").Value = 23101807
$ws.Range("C4").Value = "puwmncg75"
$ws.Range("D4").Value = "jbA6&7E!"
$ws.Range("F4").Value = "JqThvmRL"
$ws.Range("G4").Value = "TiMo"
